$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 17 into row 18/19 so the new cells pick up the
# same number formats / styles (date, time, centered text) used throughout
# the table, then overwrite with the new row's actual values.
$ws.Range("A17:C17").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H17").Copy()
$ws.Range("H18").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Row 18: fill in the new entry data (date, time in/out, hours, activity text)
$ws.Range("A18").Value2 = 43888

$ws.Range("B18").Value = 0.3125
$ws.Range("C18").Value = 0.39583333333333331

$ws.Range("F18").Value = 2

# G18 already holds a formula (=F18+G17); recalculated automatically once F18 is set
$ws.Range("G18").Formula = "=F18+G17"

$ws.Range("H18").Value = "PDF e correção de erro na tela de login e gerência de usuário"

# Row 19: running total formula continuing from G18
$ws.Range("G19").Formula = "=G18+F19"

# Update selection to match the new active cell recorded in the diff
$ws.Range("F19").Select()

$wb.Application.Calculate()
